# Update the "want to go" counts (F column) on both the "展览" (Exhibitions)
# sheet and the "全部类型" (All Types) sheet, which mirror the same data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 31
    $ws.Range("F4").Value = 112
    $ws.Range("F5").Value = 38
}
